# Apply the "AWS bash cmd line" update to bbCeny.xlsx
#  - refresh the "Last status check on" timestamp in F1
#  - record the new check results for row 8 (Benzina Albert Modrice):
#      B8 = new price, C8 = previous price (was B8),
#      D8 = delta written as a signed text string ("+0.3"),
#      E8 = the check timestamp written as plain text (no date formatting)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the status-check banner
$ws.Range("F1").Value = "Last status check on: 13.02.2022 09:15"

# 2) Update row 8's price figures
$ws.Range("B8").Value = 37.2
$ws.Range("C8").Value = 36.9

# D8: delta is now stored as literal text "+0.3" (not a number), with default/no special style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "+0.3"
$ws.Range("D8").Style = "Normal"

# E8: timestamp now stored as literal text (not a serial date), losing its date number format
$ws.Range("E8").Style = "Normal"
$ws.Range("E8").Value = "2022-02-13 09:15:16"
